# Populate IVA_CONTROL sheet with data from TRANSACCIONES (Nov 2025)
# per poblar_iva_desde_transacciones.py
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IVA_CONTROL")

# ---- VENTAS - IVA COBRADO (rows 6-20) ----
$ws.Range("A6").Value = 45970.96726364583
$ws.Range("B6").Value = 'T-41'
$ws.Range("C6").Value = 'CPF SERVICIOS RADIOLÓGICOS'
$ws.Range("D6").Value = 56.5
$ws.Range("K6").Value = 41

$ws.Range("A7").Value = 45970.96726364583
$ws.Range("B7").Value = 'T-42'
$ws.Range("C7").Value = 'ORTODEC S.A.'
$ws.Range("D7").Value = 56.5
$ws.Range("G7").Value = 'NO'
$ws.Range("J7").Value = 'EMITIDA'
$ws.Range("K7").Value = 42

$ws.Range("A8").Value = 45970.96726364583
$ws.Range("B8").Value = 'T-43'
$ws.Range("C8").Value = 'ORTODONCIA DE LA CRUZ'
$ws.Range("D8").Value = 356.5
$ws.Range("G8").Value = 'NO'
$ws.Range("J8").Value = 'EMITIDA'
$ws.Range("K8").Value = 43

$ws.Range("A9").Value = 45970.96726364583
$ws.Range("B9").Value = 'T-44'
$ws.Range("C9").Value = 'SMART WEB SERVICES'
$ws.Range("D9").Value = 1237.35
$ws.Range("G9").Value = 'NO'
$ws.Range("J9").Value = 'EMITIDA'
$ws.Range("K9").Value = 44

$ws.Range("A10").Value = 45970.96726365741
$ws.Range("B10").Value = 'T-53'
$ws.Range("C10").Value = 'GRUPO PORCINAS DE LA COSTA'
$ws.Range("D10").Value = 1171.18
$ws.Range("G10").Value = 'NO'
$ws.Range("J10").Value = 'EMITIDA'
$ws.Range("K10").Value = 53

$ws.Range("A11").Value = 45970.96726365741
$ws.Range("B11").Value = 'T-54'
$ws.Range("C11").Value = 'RODRIGUEZ ROJAS CARLOS HUMBERTO'
$ws.Range("D11").Value = 282.5
$ws.Range("G11").Value = 'NO'
$ws.Range("J11").Value = 'EMITIDA'
$ws.Range("K11").Value = 54

$ws.Range("A12").Value = 45970.96726365741
$ws.Range("B12").Value = 'T-55'
$ws.Range("C12").Value = 'VOLIO PARTNERS'
$ws.Range("D12").Value = 284.76
$ws.Range("G12").Value = 'NO'
$ws.Range("J12").Value = 'EMITIDA'
$ws.Range("K12").Value = 55

$ws.Range("A13").Value = 45970.96726365741
$ws.Range("B13").Value = 'T-56'
$ws.Range("C13").Value = 'SMART WEB SERVICES'
$ws.Range("D13").Value = 149.16
$ws.Range("G13").Value = 'NO'
$ws.Range("J13").Value = 'EMITIDA'
$ws.Range("K13").Value = 56

$ws.Range("A14").Value = 45970.96726365741
$ws.Range("B14").Value = 'T-57'
$ws.Range("C14").Value = 'GENTRA DE COSTA RICA'
$ws.Range("D14").Value = 226.0
$ws.Range("G14").Value = 'NO'
$ws.Range("J14").Value = 'EMITIDA'
$ws.Range("K14").Value = 57

$ws.Range("A15").Value = 45970.96726365741
$ws.Range("B15").Value = 'T-62'
$ws.Range("C15").Value = 'ALMACEN FISCAL DEL PACIFICO'
$ws.Range("D15").Value = 761.06
$ws.Range("G15").Value = 'NO'
$ws.Range("J15").Value = 'EMITIDA'
$ws.Range("K15").Value = 62

$ws.Range("A16").Value = 45970.97982065972
$ws.Range("B16").Value = 'T-75'
$ws.Range("C16").Value = 'OPERATION MANAGMENT OF TIERRA MAGNIFICA'
$ws.Range("D16").Value = 209.06
$ws.Range("G16").Value = 'NO'
$ws.Range("J16").Value = 'EMITIDA'
$ws.Range("K16").Value = 75

$ws.Range("A17").Value = 45970.97982065972
$ws.Range("B17").Value = 'T-76'
$ws.Range("C17").Value = 'CPF SERVICIOS RADIOLÓGICOS S.A.'
$ws.Range("D17").Value = 56.5
$ws.Range("G17").Value = 'NO'
$ws.Range("J17").Value = 'EMITIDA'
$ws.Range("K17").Value = 76

$ws.Range("A18").Value = 45970.97982065972
$ws.Range("B18").Value = 'T-77'
$ws.Range("C18").Value = 'ORTODEC S.A.'
$ws.Range("D18").Value = 56.5
$ws.Range("G18").Value = 'NO'
$ws.Range("J18").Value = 'EMITIDA'
$ws.Range("K18").Value = 77

$ws.Range("A19").Value = 45970.97982065972
$ws.Range("B19").Value = 'T-78'
$ws.Range("C19").Value = 'CEMSO'
$ws.Range("D19").Value = 333.92
$ws.Range("G19").Value = 'NO'
$ws.Range("J19").Value = 'EMITIDA'
$ws.Range("K19").Value = 78

$ws.Range("A20").Value = 45970.97982065972
$ws.Range("B20").Value = 'T-79'
$ws.Range("C20").Value = 'GRUPO ACCION COMERCIAL S.A.'
$ws.Range("D20").Value = 1689.04
$ws.Range("G20").Value = 'NO'
$ws.Range("J20").Value = 'EMITIDA'
$ws.Range("K20").Value = 79

# ---- COMPRAS - IVA PAGADO (rows 25-40) ----
$ws.Range("A25").Value = 45970.39496228009
$ws.Range("B25").Value = 'C-9'
$ws.Range("C25").Value = 'VWR INTERNATIONAL LTDA'
$ws.Range("D25").Value = 2477.876106194691
$ws.Range("I25").Value = 'N/D'
$ws.Range("K25").Value = 9

$ws.Range("A26").Value = 45970.39496228009
$ws.Range("B26").Value = 'C-10'
$ws.Range("C26").Value = 'GRUPO ACCION COMERCIAL S.A.'
$ws.Range("D26").Value = 1494.725663716814
$ws.Range("G26").Value = 'SI'
$ws.Range("I26").Value = 'N/D'
$ws.Range("J26").Value = 'PAGADA'
$ws.Range("K26").Value = 10

$ws.Range("A27").Value = 45970.39496228009
$ws.Range("B27").Value = 'C-11'
$ws.Range("C27").Value = 'ALMACEN FISCAL DEL PACIFICO ALFIPAC SOCIEDAD ANONIMA'
$ws.Range("D27").Value = 673.4955752212389
$ws.Range("G27").Value = 'SI'
$ws.Range("I27").Value = 'N/D'
$ws.Range("J27").Value = 'PAGADA'
$ws.Range("K27").Value = 11

$ws.Range("A28").Value = 45970.39496228009
$ws.Range("B28").Value = 'C-12'
$ws.Range("C28").Value = '3-102-887892 SOCIEDAD DE RESPONSABILIDAD LIMITADA'
$ws.Range("D28").Value = 612.0
$ws.Range("G28").Value = 'SI'
$ws.Range("I28").Value = 'N/D'
$ws.Range("J28").Value = 'PAGADA'
$ws.Range("K28").Value = 12

$ws.Range("A29").Value = 45970.39496228009
$ws.Range("B29").Value = 'C-13'
$ws.Range("C29").Value = 'WAIPIO SOCIEDAD ANONIMA'
$ws.Range("D29").Value = 608.2035398230089
$ws.Range("G29").Value = 'SI'
$ws.Range("I29").Value = 'N/D'
$ws.Range("J29").Value = 'PAGADA'
$ws.Range("K29").Value = 13

$ws.Range("A30").Value = 45970.39496228009
$ws.Range("B30").Value = 'C-14'
$ws.Range("C30").Value = 'CENTRO INTEGRAL DE ONCOLOGIA CIO SRL'
$ws.Range("D30").Value = 608.0088495575221
$ws.Range("G30").Value = 'SI'
$ws.Range("I30").Value = 'N/D'
$ws.Range("J30").Value = 'PAGADA'
$ws.Range("K30").Value = 14

$ws.Range("A31").Value = 45970.39496228009
$ws.Range("B31").Value = 'C-15'
$ws.Range("C31").Value = 'ORTODONCIA DE LA CRUZ'
$ws.Range("D31").Value = 437.6106194690266
$ws.Range("G31").Value = 'SI'
$ws.Range("I31").Value = 'N/D'
$ws.Range("J31").Value = 'PAGADA'
$ws.Range("K31").Value = 15

$ws.Range("A32").Value = 45970.39496228009
$ws.Range("B32").Value = 'C-16'
$ws.Range("C32").Value = 'GLOBAL AUTOMOTRIZ DE COSTA RICA GACR SOCIEDAD ANONIMA'
$ws.Range("D32").Value = 389.0353982300886
$ws.Range("G32").Value = 'SI'
$ws.Range("I32").Value = 'N/D'
$ws.Range("J32").Value = 'PAGADA'
$ws.Range("K32").Value = 16

$ws.Range("A33").Value = 45970.39496228009
$ws.Range("B33").Value = 'C-17'
$ws.Range("C33").Value = 'SOLUSA CONSOLIDATORS AND FORWARDIG, S.A.'
$ws.Range("D33").Value = 334.8230088495576
$ws.Range("G33").Value = 'SI'
$ws.Range("I33").Value = 'N/D'
$ws.Range("J33").Value = 'PAGADA'
$ws.Range("K33").Value = 17

$ws.Range("A34").Value = 45970.39496228009
$ws.Range("B34").Value = 'C-18'
$ws.Range("C34").Value = 'CEMSO'
$ws.Range("D34").Value = 295.5044247787611
$ws.Range("G34").Value = 'SI'
$ws.Range("I34").Value = 'N/D'
$ws.Range("J34").Value = 'PAGADA'
$ws.Range("K34").Value = 18

$ws.Range("A35").Value = 45970.39496228009
$ws.Range("B35").Value = 'C-19'
$ws.Range("C35").Value = 'ASOCIACION COSTARRICENSE DE AGENCIAS DE CARGA Y LOGISTICA INTERNACIONAL ACACIA'
$ws.Range("D35").Value = 295.0000000000001
$ws.Range("G35").Value = 'SI'
$ws.Range("I35").Value = 'N/D'
$ws.Range("J35").Value = 'PAGADA'
$ws.Range("K35").Value = 19

$ws.Range("A36").Value = 45970.39496228009
$ws.Range("B36").Value = 'C-20'
$ws.Range("C36").Value = 'RODRIGUEZ ROJAS CARLOS HUMBERTO'
$ws.Range("D36").Value = 250.0
$ws.Range("G36").Value = 'SI'
$ws.Range("I36").Value = 'N/D'
$ws.Range("J36").Value = 'PAGADA'
$ws.Range("K36").Value = 20

$ws.Range("A37").Value = 45970.39496228009
$ws.Range("B37").Value = 'C-21'
$ws.Range("C37").Value = 'SUPPLY NET C.R.W.H SOCIEDAD ANONIMA'
$ws.Range("D37").Value = 245.0000000000001
$ws.Range("G37").Value = 'SI'
$ws.Range("I37").Value = 'N/D'
$ws.Range("J37").Value = 'PAGADA'
$ws.Range("K37").Value = 21

$ws.Range("A38").Value = 45970.39496228009
$ws.Range("B38").Value = 'C-22'
$ws.Range("C38").Value = 'OPERATION MANAGMENT OF TIERRA MAGNIFICA'
$ws.Range("D38").Value = 185.0088495575222
$ws.Range("G38").Value = 'SI'
$ws.Range("I38").Value = 'N/D'
$ws.Range("J38").Value = 'PAGADA'
$ws.Range("K38").Value = 22

$ws.Range("A39").Value = 45970.39496228009
$ws.Range("B39").Value = 'C-23'
$ws.Range("C39").Value = 'GENTRA DE COSTA RICA SOCIEDAD ANONIMA'
$ws.Range("D39").Value = 162.5044247787611
$ws.Range("G39").Value = 'SI'
$ws.Range("I39").Value = 'N/D'
$ws.Range("J39").Value = 'PAGADA'
$ws.Range("K39").Value = 23

$ws.Range("A40").Value = 45970.39496228009
$ws.Range("B40").Value = 'C-24'
$ws.Range("C40").Value = 'SEVILLA NAVARRO EDGAR'
$ws.Range("D40").Value = 150.0
$ws.Range("G40").Value = 'SI'
$ws.Range("I40").Value = 'N/D'
$ws.Range("J40").Value = 'PAGADA'
$ws.Range("K40").Value = 24
